# Change the table style of the three balance-sheet tables (slides 14-16)
# from the custom "Table_0" style to the built-in PowerPoint table style
# {F1DD60CA-895D-4C30-966F-70757090B508}.

$p = $ppt.ActivePresentation

$targetSlides = @(14, 15, 16)
$newStyleId = "{F1DD60CA-895D-4C30-966F-70757090B508}"

foreach ($slideIdx in $targetSlides) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newStyleId, $true)
        }
    }
}
